$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a brand-new row 12 for the split-off "Dashboard view" task first,
#     so later formatting tweaks on row 11 don't leak into the new row ---
$ws.Rows("12:12").Insert()

# --- Fill in the new row 12 content (task split off of "Dashboard") ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Front-end: Dashboard view"
$ws.Range("C12").Value = "LMS v3"
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 44231
$ws.Range("F12").Value = 44237

# --- Row 11: rename task, add Actual Finish date + Comment ---
$ws.Range("B11").Value = "Front-end: Dashboard main page"
$ws.Range("G10").Copy()
$ws.Range("G11").PasteSpecial(-4122)   # xlPasteFormats, reuse the date-format style
$ws.Range("G11").Value = 44230
$ws.Range("I11").Value = "Dashboard as a whole took longer than expected, need more time to complete the view page"

# --- Remaining rows shift down by one; update their S/N + dates/duration ---
# (old row 12 "Front-end: Integrate Threejs scene" -> now row 13)
$ws.Range("A13").Value = 12
$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 44238
$ws.Range("F13").Value = 44246

# (old row 13 "Front-end: User management/Settings" -> now row 14)
$ws.Range("A14").Value = 13
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 44249
$ws.Range("F14").Value = 44251

# (old row 14 "Front-end: Dark mode toggle" -> now row 15)
$ws.Range("A15").Value = 14
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 44252
$ws.Range("F15").Value = 44253

# (old row 15 "Documentation for existing/planned features" -> now row 16, content unchanged)
$ws.Range("A16").Value = 15

# --- Update the last-selected cell shown in the sheet view ---
$ws.Range("H20").Select()
